$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format column D as Text so numeric-looking strings
# (e.g. "1.001", "7.080") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.901.85"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "1.649.41"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "308.54"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "0.3889"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").Value = "0.3826"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Value = "52.05"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "1.351"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "0.08419"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "23.84"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "7.080"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "7.949"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "0.00001314"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "1.647.02"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "94.66"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "0.06964"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "19.68"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "6.933"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "13.72"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").Value = "23.895.07"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "2.959"
$ws.Range("D27").Value = "22.09"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").Value = "151.56"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("D29").Value = "5.395"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").Value = "139.01"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "2.517"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "1.828.35"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "1.042"
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("D35").Value = "0.08052"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "10.96"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("D38").Value = "6.654"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").Value = "0.2679"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "0.09091"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "0.7619"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "13.44"
$ws.Range("E42").Value = "  -3.28%  "
$ws.Range("D43").Value = "1.427"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "16.40"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").Value = "0.7004"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "2.466"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "4.074"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "0.08298"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "1.213"
$ws.Range("E51").Value = "  -2.27%  "

# Restore the original (default/Normal) style on column D so the
# cell style index is unchanged from before the edit.
$ws.Range("D2:D51").Style = "Normal"

